# Generate Report for Handback
# Updates timestamps / status strings on the handback-status workbook.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-20 09:52:13"
$wsOverview.Range("G3").Value = "2016-10-20 09:52:13"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-10-20 09:52:01"
$wsZhCn.Range("H3").Value = "2016-10-20 09:52:01"
# Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-10-20 09:52:45"
$wsZhCn.Range("K3").Value = "2016-10-20 09:52:45"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Correspond Handoff Datetime column (H)
$wsDeDe.Range("H2").Value = "2016-10-20 09:52:13"
$wsDeDe.Range("H3").Value = "2016-10-20 09:52:13"
# Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-10-20 09:53:03"
$wsDeDe.Range("K3").Value = "2016-10-20 09:53:03"
